# Updating glossary and use cases
$d = $word.ActiveDocument
$apos = [char]8217

# 1) Title: "Use cases for module Customer:" -> "Use cases for module Customer Service:"
$d.Content.Find.Execute(
    "Use cases for module Customer:", $true, $false, $false, $false, $false,
    $true, 0, $false, "Use cases for module Customer Service:", 1) | Out-Null

# 2) Table 1, Row 1 (Positive) - adding item into the actor's own basket
$t1 = $d.Tables.Item(1)
$r = $t1.Cell(1, 2).Range
$r.Find.Execute(
    "try to add an item from an existing shop into the shop basket of a logged in user.",
    $true, $false, $false, $false, $false,
    $true, 0, $false,
    "Try to add an item from an existing shop into the shop basket of the actor.", 1) | Out-Null

# 3) Table 1, Row 2 (Negative) - adding item into another user's basket
$t1b = $d.Tables.Item(1)
$r = $t1b.Cell(2, 2).Range
$r.Find.Execute(
    "try to add an item from an existing shop into the shop basket of a not logged in user.",
    $true, $false, $false, $false, $false,
    $true, 0, $false,
    "Try to add an item from an existing shop into the shop basket of another (not the actor) user.", 1) | Out-Null

# 4) Table 2, Row 1 (Positive) - removing an item from the actor's own basket
$t2 = $d.Tables.Item(2)
$r = $t2.Cell(1, 2).Range
$r.Find.Execute(
    "a logged in user" + $apos + "s basket.",
    $true, $false, $false, $false, $false,
    $true, 0, $false,
    "the actor" + $apos + "s basket", 1) | Out-Null

# 5) Table 2, Row 2 (Negative) - removing an item NOT from the actor's own basket
$t2b = $d.Tables.Item(2)
$r = $t2b.Cell(2, 2).Range
$r.Find.Execute(
    "an item from a logged in user" + $apos + "s basket.",
    $true, $false, $false, $false, $false,
    $true, 0, $false,
    "an item not from the actor" + $apos + "s basket", 1) | Out-Null

# 6) Table 3, Row 1 (Positive) - deleting the customer account
$t3 = $d.Tables.Item(3)
$r = $t3.Cell(1, 2).Range
$r.Find.Execute(
    "A logged in Admin user (or regular user if it plans to remove own account) deletes the customer account",
    $true, $false, $false, $false, $false,
    $true, 0, $false,
    "An Admin user (or regular user if removing own account) deletes the customer account", 1) | Out-Null

# 7) Table 4, Row 1 (Positive) - updating the customer account
$t4 = $d.Tables.Item(4)
$r = $t4.Cell(1, 2).Range
$r.Find.Execute(
    "A logged in Admin user (or regular user if it plans to update own account) updates the customer account",
    $true, $false, $false, $false, $false,
    $true, 0, $false,
    "An Admin user (or regular user if updates own account) updates the customer account", 1) | Out-Null
